$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C9 previously held "TI, SO, PU, SC" - remove the nonexistent "TI" column
$ws.Range("C9").Value = "SO, PU, SC"

# Update the selection to reflect the now-relevant cell (C9)
$ws.Range("C9").Select()
